$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9690061211585999
$ws.Range("B1").Value = 1.909066796302795
$ws.Range("C1").Value = 5.027409076690674
$ws.Range("D1").Value = 1.793805360794067
$ws.Range("E1").Value = 0.6933905482292175
